$wb = $excel.ActiveWorkbook

# --- Cabling sheet: change phase voltage/current usage and phase values ---
$cabling = $wb.Worksheets.Item("Cabling")

# C5:C19 -> AMPS (was VOLTS)
$cabling.Range("C5:C19").Value = "AMPS"

# G5:G7 -> NONE (was A/AB, B/BC, C/CA)
$cabling.Range("G5:G7").Value = "NONE"

# Column F widened to fit the new channel-number labels
$cabling.Columns.Item(6).ColumnWidth = 13.0

# Selection for the Cabling sheet moves to E17:J19
$cabling.Range("E17:J19").Select()

# --- DSPChannelMap sheet: rename channel-number headers ---
$dsp = $wb.Worksheets.Item("DSPChannelMap")
$dsp.Range("B1").Value = "Channel_Number_DSP1"
$dsp.Range("D1").Value = "Channel_Number_DSP2"

# Widen columns B/D to fit the longer header labels
$dsp.Columns.Item(2).ColumnWidth = 21.6
$dsp.Columns.Item(4).ColumnWidth = 21.6

$dsp.Range("B1").Select()

# --- Make Cabling the active sheet/tab ---
$cabling.Activate()
